$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 10:22"

# --- Reorder countries whose ranking swapped (Polonia overtook Japon; Filipinas overtook Malasia) ---
$ws.Range("A30").Value = "Polonia"
$ws.Range("A31").Value = "Japon"

$ws.Range("A37").Value = "Filipinas"
$ws.Range("A38").Value = "Malasia"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 19: Austria
$ws.Range("B19").Value = 13667
$ws.Range("C19").Value = 107
$ws.Range("E19").Value = 7284

# Row 20: Rusia
$ws.Range("B20").Value = 13584
$ws.Range("C20").Value = 1667
$ws.Range("D20").Value = 1045
$ws.Range("E20").Value = 12433
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 106

# Row 30: Polonia (moved up)
$ws.Range("B30").Value = 6088
$ws.Range("C30").Value = 133
$ws.Range("D30").Value = 375
$ws.Range("E30").Value = 5518
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 14
$ws.Range("H30").Value = 195

# Row 31: Japon (moved down)
$ws.Range("B31").Value = 6005
$ws.Range("D31").Value = 685
$ws.Range("E31").Value = 5221
$ws.Range("F31").Value = 109
$ws.Range("H31").Value = 99

# Row 37: Filipinas (moved up)
$ws.Range("B37").Value = 4428
$ws.Range("C37").Value = 233
$ws.Range("D37").Value = 157
$ws.Range("E37").Value = 4024
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 26
$ws.Range("H37").Value = 247

# Row 38: Malasia (moved down)
$ws.Range("B38").Value = 4346
$ws.Range("D38").Value = 1830
$ws.Range("E38").Value = 2446
$ws.Range("F38").Value = 69
$ws.Range("H38").Value = 70

# Row 50: Ucrania
$ws.Range("F50").Value = 45

# Row 74: Bosnia y Herzegovina
$ws.Range("B74").Value = 917
$ws.Range("C74").Value = 16
$ws.Range("D74").Value = 137
$ws.Range("E74").Value = 743
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 37

# Row 80: Tunez
$ws.Range("D80").Value = 43
$ws.Range("E80").Value = 603
$ws.Range("F80").Value = 85

# Row 83: Letonia
$ws.Range("B83").Value = 630
$ws.Range("C83").Value = 18
$ws.Range("E83").Value = 611
$ws.Range("F83").Value = 2

# Row 102: Malta
$ws.Range("E102").Value = 331
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 3

# Row 110: Montenegro
$ws.Range("D110").Value = 5
$ws.Range("E110").Value = 255
